$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended at the bottom of the sheet (rows 5-7), duplicating the
# existing match rows (re-ordered), as described by the diff.
# Note: the batsman name in the source data ends with a non-breaking space
# (U+00A0), matching the existing rows 2-4 in this sheet - use that exact
# character (not a regular space) so the duplicated rows are byte-faithful.
$nbsp = [char]0x00A0
$batsman = "Josh Philippe" + $nbsp

$rows = @(
    @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Royal Challengers Bangalore", "Mumbai Indians", $batsman, "33", "24", "4", "1", "137.50"),
    @(" Abu Dhabi", " November 02 2020", "Capitals won by 6 wickets (with 6 balls remaining)", "Royal Challengers Bangalore", "Delhi Capitals", $batsman, "12", "17", "1", "0", "70.58"),
    @(" Sharjah", " October 31 2020", "Sunrisers won by 5 wickets (with 35 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", $batsman, "32", "31", "4", "0", "103.22")
)

$startRow = 5
$endRow = $startRow + $rows.Length - 1
$targetRange = $ws.Range("A" + $startRow + ":K" + $endRow)

# Every value in this sheet (including the numeric-looking ones such as
# totalRuns/sr) is stored as text, matching the existing rows. Pre-set the
# destination range's number format to Text so Excel doesn't coerce
# numeric-looking strings (e.g. "33", "137.50") into real numbers when
# they're written below.
$targetRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

# Drop the temporary Text number format now that the values are safely
# stored as text, so the new cells end up with the sheet's default style
# (matching rows 1-4, which carry no explicit cell format).
$targetRange.ClearFormats()
